# Update "想去人数" (want-to-go count) figures in the F column for a handful
# of rows on the "展览" and "全部类型" worksheets, reflecting the site's
# refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows keyed by row number -> old/new F value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4247
$ws1.Range("F3").Value = 2416
$ws1.Range("F7").Value = 49
$ws1.Range("F10").Value = 126
$ws1.Range("F12").Value = 1583
$ws1.Range("F14").Value = 3255

# Sheet "全部类型" (same events, different row numbers because this sheet
# aggregates all categories)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4247
$ws4.Range("F3").Value = 2416
$ws4.Range("F8").Value = 49
$ws4.Range("F12").Value = 126
$ws4.Range("F16").Value = 1583
$ws4.Range("F18").Value = 3255
